# Rejected-data check: the old "Min_WQ_Crit vs DEQ_Pref_Method" comparison column (H)
# is no longer needed on this sheet, so remove it entirely (header, data and the
# conditional formatting that lived on it), rename the sheet to the generic
# "Sheet1", and refresh the AutoFilter / FilterDatabase range to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename sheet -------------------------------------------------------
$ws.Name = "Sheet1"

# --- drop column H (and any stray columns out to V) ---------------------
$ws.Range("H1:V191").Clear()

# --- column widths (manually re-tuned after the column removal) ---------
$ws.Columns("A").ColumnWidth = 34.5703125
$ws.Columns("B").ColumnWidth = 9.140625
$ws.Columns("C").ColumnWidth = 37
$ws.Columns("D").ColumnWidth = 16.7109375
$ws.Columns("E").ColumnWidth = 12.140625
$ws.Columns("F").ColumnWidth = 9.140625
$ws.Columns("G").ColumnWidth = 13.5703125

# --- freeze the header row and scroll down a bit -------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E1").Select()

# --- autofilter on the new (narrower) range ------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:G191").AutoFilter()
